# "frei wählbare Quarzfrequenz" now fully achieved: IST (D24) goes from 0 to 5.
# D34 (Gesamt / SUM(D2:D33)) and D35 (D34/60) recalc automatically as formulas.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D24").Value = 5

# Update the view's active cell/selection to D18 (matches the saved view state).
$ws.Range("D18").Select()
